# PNAD 2009 - "correção nos dados e inicio da analise"
# Remove the "sexo" sub-header row (row 5) and the "cor ou raça" sub-header
# row (originally row 8, now row 7 after the first deletion) so that the
# data rows (homens/mulheres/branca/preta ou parda) move up and directly
# follow the "total (1)" row, each keeping its own numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(5).Delete()   # "sexo" header row (was row 5)
$ws.Rows.Item(7).Delete()   # "cor ou raça" header row (was row 8, now row 7)
